$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 22
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 0.727
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 0.456
$ws.Range("M2").Value = 17.5
$ws.Range("N2").Value = 0.404
$ws.Range("O2").Value = 16.6
$ws.Range("P2").Value = 22.7
$ws.Range("Q2").Value = 0.732
$ws.Range("R2").Value = 10.2
$ws.Range("S2").Value = 32.1
$ws.Range("T2").Value = 42.3
$ws.Range("V2").Value = 14
$ws.Range("W2").Value = 8.199999999999999
$ws.Range("AA2").Value = 19.8
$ws.Range("AB2").Value = 97.59999999999999
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 8
$ws.Range("AF2").Value = 2
$ws.Range("AG2").Value = 3
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 8
$ws.Range("AJ2").Value = 12
$ws.Range("AL2").Value = 7
$ws.Range("AO2").Value = 15
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 22
$ws.Range("AR2").Value = 21
$ws.Range("AS2").Value = 8
$ws.Range("AT2").Value = 14
$ws.Range("AW2").Value = 13
$ws.Range("AY2").Value = 12
$ws.Range("BA2").Value = 17
$ws.Range("BB2").Value = 8
$ws.Range("BC2").Value = 4
$ws.Range("BF2").Value = "2012-02-02"

# Row 3
$ws.Range("AD3").Value = 18
$ws.Range("AE3").Value = 16
$ws.Range("AF3").Value = 15
$ws.Range("AG3").Value = 16
$ws.Range("AP3").Value = 20
$ws.Range("AS3").Value = 23
$ws.Range("AY3").Value = 21
$ws.Range("BA3").Value = 18
$ws.Range("BF3").Value = "2012-02-02"

# Row 4
$ws.Range("AD4").Value = 3
$ws.Range("AH4").Value = 19
$ws.Range("AJ4").Value = 8
$ws.Range("AP4").Value = 25
$ws.Range("AQ4").Value = 21
$ws.Range("AS4").Value = 25
$ws.Range("AT4").Value = 26
$ws.Range("AV4").Value = 13
$ws.Range("AZ4").Value = 11
$ws.Range("BF4").Value = "2012-02-02"

# Row 5
$ws.Range("D5").Value = 24
$ws.Range("E5").Value = 18
$ws.Range("G5").Value = 0.75
$ws.Range("I5").Value = 37.3
$ws.Range("J5").Value = 81.5
$ws.Range("K5").Value = 0.458
$ws.Range("M5").Value = 14.9
$ws.Range("N5").Value = 0.364
$ws.Range("O5").Value = 16.3
$ws.Range("P5").Value = 22.2
$ws.Range("R5").Value = 13.2
$ws.Range("S5").Value = 31.8
$ws.Range("T5").Value = 44.9
$ws.Range("U5").Value = 22.3
$ws.Range("W5").Value = 7.2
$ws.Range("X5").Value = 6.1
$ws.Range("Y5").Value = 4.8
$ws.Range("Z5").Value = 17.3
$ws.Range("AA5").Value = 18
$ws.Range("AB5").Value = 96.3
$ws.Range("AC5").Value = 8.4
$ws.Range("AI5").Value = 6
$ws.Range("AJ5").Value = 10
$ws.Range("AK5").Value = 8
$ws.Range("AL5").Value = 23
$ws.Range("AM5").Value = 21
$ws.Range("AN5").Value = 9
$ws.Range("AO5").Value = 18
$ws.Range("AP5").Value = 17
$ws.Range("AQ5").Value = 20
$ws.Range("AS5").Value = 12
$ws.Range("AT5").Value = 1
$ws.Range("AX5").Value = 4
$ws.Range("AY5").Value = 10
$ws.Range("BB5").Value = 12
$ws.Range("BF5").Value = "2012-02-02"

# Row 6
$ws.Range("AD6").Value = 27
$ws.Range("AH6").Value = 12
$ws.Range("AJ6").Value = 19
$ws.Range("AL6").Value = 8
$ws.Range("AO6").Value = 13
$ws.Range("AS6").Value = 26
$ws.Range("AT6").Value = 18
$ws.Range("AV6").Value = 29
$ws.Range("BB6").Value = 18
$ws.Range("BF6").Value = "2012-02-02"

# Row 7
$ws.Range("AD7").Value = 3
$ws.Range("AE7").Value = 7
$ws.Range("AF7").Value = 10
$ws.Range("AG7").Value = 9
$ws.Range("AH7").Value = 19
$ws.Range("AL7").Value = 9
$ws.Range("AQ7").Value = 18
$ws.Range("AR7").Value = 20
$ws.Range("AV7").Value = 6
$ws.Range("BC7").Value = 8
$ws.Range("BF7").Value = "2012-02-02"

# Row 8
$ws.Range("D8").Value = 21
$ws.Range("E8").Value = 14
$ws.Range("G8").Value = 0.667
$ws.Range("H8").Value = 49
$ws.Range("J8").Value = 80.3
$ws.Range("K8").Value = 0.481
$ws.Range("L8").Value = 6.9
$ws.Range("N8").Value = 0.338
$ws.Range("P8").Value = 28.6
$ws.Range("Q8").Value = 0.744
$ws.Range("R8").Value = 9.4
$ws.Range("S8").Value = 33.1
$ws.Range("U8").Value = 24
$ws.Range("W8").Value = 9.699999999999999
$ws.Range("Y8").Value = 6.3
$ws.Range("Z8").Value = 18.5
$ws.Range("AA8").Value = 24
$ws.Range("AB8").Value = 105.4
$ws.Range("AC8").Value = 6.2
$ws.Range("AD8").Value = 18
$ws.Range("AE8").Value = 7
$ws.Range("AF8").Value = 8
$ws.Range("AG8").Value = 8
$ws.Range("AJ8").Value = 18
$ws.Range("AL8").Value = 11
$ws.Range("AN8").Value = 16
$ws.Range("AQ8").Value = 16
$ws.Range("AV8").Value = 27
$ws.Range("AW8").Value = 2
$ws.Range("BC8").Value = 6
$ws.Range("BF8").Value = "2012-02-02"

# Row 9
$ws.Range("AD9").Value = 1
$ws.Range("AK9").Value = 24
$ws.Range("AN9").Value = 14
$ws.Range("AP9").Value = 24
$ws.Range("AQ9").Value = 7
$ws.Range("AR9").Value = 13
$ws.Range("AV9").Value = 26
$ws.Range("BC9").Value = 28
$ws.Range("BF9").Value = "2012-02-02"

# Row 10
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 7
$ws.Range("G10").Value = 0.368
$ws.Range("I10").Value = 37
$ws.Range("J10").Value = 80.90000000000001
$ws.Range("K10").Value = 0.457
$ws.Range("M10").Value = 19.9
$ws.Range("N10").Value = 0.369
$ws.Range("O10").Value = 14.8
$ws.Range("P10").Value = 20.5
$ws.Range("Q10").Value = 0.722
$ws.Range("R10").Value = 9.800000000000001
$ws.Range("S10").Value = 30
$ws.Range("T10").Value = 39.8
$ws.Range("V10").Value = 15.2
$ws.Range("W10").Value = 8.9
$ws.Range("AA10").Value = 17.6
$ws.Range("AB10").Value = 96.2
$ws.Range("AC10").Value = -3.1
$ws.Range("AD10").Value = 29
$ws.Range("AE10").Value = 24
$ws.Range("AG10").Value = 23
$ws.Range("AI10").Value = 7
$ws.Range("AJ10").Value = 14
$ws.Range("AL10").Value = 6
$ws.Range("AM10").Value = 12
$ws.Range("AN10").Value = 8
$ws.Range("AO10").Value = 27
$ws.Range("AP10").Value = 21
$ws.Range("AQ10").Value = 26
$ws.Range("AR10").Value = 27
$ws.Range("AS10").Value = 22
$ws.Range("AT10").Value = 27
$ws.Range("AU10").Value = 3
$ws.Range("AV10").Value = 14
$ws.Range("AW10").Value = 7
$ws.Range("BB10").Value = 13
$ws.Range("BF10").Value = "2012-02-02"

# Row 11
$ws.Range("AD11").Value = 8
$ws.Range("AF11").Value = 15
$ws.Range("AL11").Value = 10
$ws.Range("AO11").Value = 26
$ws.Range("AV11").Value = 11
$ws.Range("BC11").Value = 15
$ws.Range("BF11").Value = "2012-02-02"

# Row 12
$ws.Range("AD12").Value = 18
$ws.Range("AG12").Value = 6
$ws.Range("AJ12").Value = 16
$ws.Range("AL12").Value = 22
$ws.Range("AN12").Value = 7
$ws.Range("AO12").Value = 5
$ws.Range("AT12").Value = 2
$ws.Range("BA12").Value = 9
$ws.Range("BF12").Value = "2012-02-02"

# Row 13
$ws.Range("D13").Value = 19
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 0.6840000000000001
$ws.Range("I13").Value = 36.8
$ws.Range("J13").Value = 80.09999999999999
$ws.Range("L13").Value = 7.7
$ws.Range("M13").Value = 21.3
$ws.Range("N13").Value = 0.363
$ws.Range("O13").Value = 17.6
$ws.Range("P13").Value = 26.4
$ws.Range("Q13").Value = 0.665
$ws.Range("R13").Value = 10.9
$ws.Range("S13").Value = 29.9
$ws.Range("T13").Value = 40.8
$ws.Range("U13").Value = 21.6
$ws.Range("V13").Value = 14
$ws.Range("X13").Value = 5.6
$ws.Range("Y13").Value = 4
$ws.Range("AA13").Value = 23.1
$ws.Range("AB13").Value = 99
$ws.Range("AC13").Value = 2
$ws.Range("AD13").Value = 29
$ws.Range("AF13").Value = 2
$ws.Range("AG13").Value = 7
$ws.Range("AI13").Value = 9
$ws.Range("AJ13").Value = 21
$ws.Range("AK13").Value = 7
$ws.Range("AL13").Value = 4
$ws.Range("AO13").Value = 10
$ws.Range("AP13").Value = 6
$ws.Range("AR13").Value = 17
$ws.Range("AS13").Value = 24
$ws.Range("AT13").Value = 24
$ws.Range("AU13").Value = 11
$ws.Range("AW13").Value = 14
$ws.Range("BC13").Value = 13
$ws.Range("BF13").Value = "2012-02-02"

# Row 14
$ws.Range("AD14").Value = 8
$ws.Range("AF14").Value = 10
$ws.Range("AG14").Value = 12
$ws.Range("AH14").Value = 15
$ws.Range("AN14").Value = 27
$ws.Range("AO14").Value = 11
$ws.Range("AQ14").Value = 14
$ws.Range("AR14").Value = 14
$ws.Range("AV14").Value = 12
$ws.Range("BB14").Value = 19
$ws.Range("BF14").Value = "2012-02-02"

# Row 15
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = 11
$ws.Range("G15").Value = 0.524
$ws.Range("I15").Value = 36.5
$ws.Range("J15").Value = 81.90000000000001
$ws.Range("K15").Value = 0.445
$ws.Range("L15").Value = 3.6
$ws.Range("M15").Value = 11
$ws.Range("N15").Value = 0.323
$ws.Range("O15").Value = 17
$ws.Range("P15").Value = 22.9
$ws.Range("Q15").Value = 0.744
$ws.Range("R15").Value = 11.7
$ws.Range("S15").Value = 30.3
$ws.Range("T15").Value = 42
$ws.Range("U15").Value = 18.9
$ws.Range("V15").Value = 15.3
$ws.Range("W15").Value = 10.5
$ws.Range("Z15").Value = 19.8
$ws.Range("AA15").Value = 20
$ws.Range("AB15").Value = 93.59999999999999
$ws.Range("AC15").Value = 1.3
$ws.Range("AD15").Value = 18
$ws.Range("AE15").Value = 16
$ws.Range("AF15").Value = 15
$ws.Range("AG15").Value = 16
$ws.Range("AH15").Value = 13
$ws.Range("AI15").Value = 12
$ws.Range("AO15").Value = 12
$ws.Range("AQ15").Value = 15
$ws.Range("AR15").Value = 9
$ws.Range("AT15").Value = 15
$ws.Range("AV15").Value = 16
$ws.Range("AX15").Value = 19
$ws.Range("AY15").Value = 21
$ws.Range("AZ15").Value = 13
$ws.Range("BB15").Value = 20
$ws.Range("BC15").Value = 14
$ws.Range("BF15").Value = "2012-02-02"

# Row 16
$ws.Range("AD16").Value = 8
$ws.Range("AM16").Value = 22
$ws.Range("AP16").Value = 2
$ws.Range("AR16").Value = 24
$ws.Range("AS16").Value = 9
$ws.Range("AT16").Value = 15
$ws.Range("AW16").Value = 6
$ws.Range("BF16").Value = "2012-02-02"

# Row 17
$ws.Range("AD17").Value = 18
$ws.Range("AQ17").Value = 5
$ws.Range("AR17").Value = 9
$ws.Range("AT17").Value = 25
$ws.Range("AW17").Value = 12
$ws.Range("AY17").Value = 18
$ws.Range("BB17").Value = 14
$ws.Range("BF17").Value = "2012-02-02"

# Row 18
$ws.Range("AD18").Value = 8
$ws.Range("AO18").Value = 4
$ws.Range("AY18").Value = 26
$ws.Range("AZ18").Value = 12
$ws.Range("BB18").Value = 11
$ws.Range("BC18").Value = 16
$ws.Range("BF18").Value = "2012-02-02"

# Row 19
$ws.Range("AD19").Value = 3
$ws.Range("AH19").Value = 19
$ws.Range("AK19").Value = 25
$ws.Range("AY19").Value = 20
$ws.Range("BF19").Value = "2012-02-02"

# Row 20
$ws.Range("D20").Value = 22
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 0.182
$ws.Range("J20").Value = 78.5
$ws.Range("M20").Value = 12.7
$ws.Range("N20").Value = 0.293
$ws.Range("O20").Value = 15
$ws.Range("P20").Value = 20.2
$ws.Range("R20").Value = 11.5
$ws.Range("S20").Value = 30
$ws.Range("T20").Value = 41.5
$ws.Range("U20").Value = 19.4
$ws.Range("W20").Value = 7.6
$ws.Range("Y20").Value = 5.8
$ws.Range("Z20").Value = 20.8
$ws.Range("AA20").Value = 18.4
$ws.Range("AB20").Value = 88.09999999999999
$ws.Range("AC20").Value = -5.2
$ws.Range("AD20").Value = 8
$ws.Range("AF20").Value = 27
$ws.Range("AG20").Value = 27
$ws.Range("AH20").Value = 15
$ws.Range("AN20").Value = 28
$ws.Range("AO20").Value = 24
$ws.Range("AP20").Value = 23
$ws.Range("AQ20").Value = 17
$ws.Range("AR20").Value = 12
$ws.Range("AV20").Value = 15
$ws.Range("AY20").Value = 26
$ws.Range("BA20").Value = 24
$ws.Range("BF20").Value = "2012-02-02"

# Row 21
$ws.Range("D21").Value = 21
$ws.Range("F21").Value = 13
$ws.Range("G21").Value = 0.381
$ws.Range("I21").Value = 33.8
$ws.Range("J21").Value = 80.09999999999999
$ws.Range("K21").Value = 0.421
$ws.Range("L21").Value = 7.4
$ws.Range("M21").Value = 23.5
$ws.Range("N21").Value = 0.314
$ws.Range("O21").Value = 19.7
$ws.Range("P21").Value = 25.2
$ws.Range("Q21").Value = 0.781
$ws.Range("R21").Value = 10.8
$ws.Range("S21").Value = 31.1
$ws.Range("T21").Value = 41.9
$ws.Range("V21").Value = 17
$ws.Range("W21").Value = 9.6
$ws.Range("X21").Value = 4.3
$ws.Range("Y21").Value = 5
$ws.Range("AA21").Value = 22
$ws.Range("AB21").Value = 94.59999999999999
$ws.Range("AC21").Value = -0.5
$ws.Range("AD21").Value = 18
$ws.Range("AF21").Value = 22
$ws.Range("AG21").Value = 21
$ws.Range("AJ21").Value = 20
$ws.Range("AK21").Value = 26
$ws.Range("AL21").Value = 5
$ws.Range("AN21").Value = 24
$ws.Range("AP21").Value = 9
$ws.Range("AQ21").Value = 6
$ws.Range("AR21").Value = 18
$ws.Range("AS21").Value = 16
$ws.Range("AT21").Value = 17
$ws.Range("AU21").Value = 23
$ws.Range("AW21").Value = 3
$ws.Range("AY21").Value = 13
$ws.Range("BA21").Value = 6
$ws.Range("BB21").Value = 15
$ws.Range("BF21").Value = "2012-02-02"

# Row 22
$ws.Range("AD22").Value = 18
$ws.Range("AP22").Value = 3
$ws.Range("AY22").Value = 11
$ws.Range("BF22").Value = "2012-02-02"

# Row 23
$ws.Range("AD23").Value = 8
$ws.Range("AF23").Value = 10
$ws.Range("AG23").Value = 12
$ws.Range("AH23").Value = 15
$ws.Range("AP23").Value = 7
$ws.Range("AS23").Value = 13
$ws.Range("BC23").Value = 18
$ws.Range("BF23").Value = "2012-02-02"

# Row 24
$ws.Range("AD24").Value = 8
$ws.Range("AH24").Value = 8
$ws.Range("AN24").Value = 6
$ws.Range("AO24").Value = 29
$ws.Range("AQ24").Value = 25
$ws.Range("BB24").Value = 9
$ws.Range("BF24").Value = "2012-02-02"

# Row 25
$ws.Range("AD25").Value = 18
$ws.Range("AG25").Value = 21
$ws.Range("AJ25").Value = 17
$ws.Range("AN25").Value = 15
$ws.Range("AQ25").Value = 8
$ws.Range("AR25").Value = 25
$ws.Range("AS25").Value = 18
$ws.Range("AZ25").Value = 10
$ws.Range("BF25").Value = "2012-02-02"

# Row 26
$ws.Range("D26").Value = 22
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 0.591
$ws.Range("I26").Value = 36.6
$ws.Range("J26").Value = 83
$ws.Range("L26").Value = 5.9
$ws.Range("N26").Value = 0.314
$ws.Range("O26").Value = 18.5
$ws.Range("P26").Value = 23.6
$ws.Range("Q26").Value = 0.785
$ws.Range("R26").Value = 11.2
$ws.Range("S26").Value = 32.4
$ws.Range("U26").Value = 21.6
$ws.Range("W26").Value = 8.800000000000001
$ws.Range("Z26").Value = 20.5
$ws.Range("AA26").Value = 21.5
$ws.Range("AB26").Value = 97.59999999999999
$ws.Range("AC26").Value = 6.3
$ws.Range("AD26").Value = 8
$ws.Range("AF26").Value = 10
$ws.Range("AG26").Value = 12
$ws.Range("AH26").Value = 15
$ws.Range("AI26").Value = 11
$ws.Range("AJ26").Value = 5
$ws.Range("AN26").Value = 23
$ws.Range("AO26").Value = 8
$ws.Range("AQ26").Value = 4
$ws.Range("AU26").Value = 10
$ws.Range("AW26").Value = 8
$ws.Range("AY26").Value = 16
$ws.Range("AZ26").Value = 18
$ws.Range("BA26").Value = 8
$ws.Range("BB26").Value = 7
$ws.Range("BF26").Value = "2012-02-02"

# Row 27
$ws.Range("D27").Value = 21
$ws.Range("E27").Value = 6
$ws.Range("G27").Value = 0.286
$ws.Range("J27").Value = 83.09999999999999
$ws.Range("K27").Value = 0.4
$ws.Range("M27").Value = 20
$ws.Range("O27").Value = 18.3
$ws.Range("P27").Value = 25
$ws.Range("Q27").Value = 0.731
$ws.Range("R27").Value = 14
$ws.Range("S27").Value = 29.1
$ws.Range("T27").Value = 43.2
$ws.Range("Z27").Value = 19.3
$ws.Range("AB27").Value = 90.59999999999999
$ws.Range("AC27").Value = -11.3
$ws.Range("AD27").Value = 18
$ws.Range("AE27").Value = 26
$ws.Range("AG27").Value = 26
$ws.Range("AJ27").Value = 4
$ws.Range("AM27").Value = 11
$ws.Range("AO27").Value = 9
$ws.Range("AP27").Value = 10
$ws.Range("AQ27").Value = 23
$ws.Range("AY27").Value = 25
$ws.Range("BA27").Value = 10
$ws.Range("BC27").Value = 29
$ws.Range("BF27").Value = "2012-02-02"

# Row 28
$ws.Range("D28").Value = 23
$ws.Range("E28").Value = 14
$ws.Range("G28").Value = 0.609
$ws.Range("H28").Value = 48.7
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 82.5
$ws.Range("L28").Value = 7.8
$ws.Range("M28").Value = 20.2
$ws.Range("N28").Value = 0.386
$ws.Range("O28").Value = 14.1
$ws.Range("P28").Value = 19.8
$ws.Range("Q28").Value = 0.714
$ws.Range("U28").Value = 22.8
$ws.Range("V28").Value = 13.9
$ws.Range("W28").Value = 7.1
$ws.Range("X28").Value = 4.1
$ws.Range("Z28").Value = 17
$ws.Range("AA28").Value = 18.7
$ws.Range("AB28").Value = 97.90000000000001
$ws.Range("AC28").Value = 3.5
$ws.Range("AD28").Value = 3
$ws.Range("AE28").Value = 7
$ws.Range("AF28").Value = 10
$ws.Range("AJ28").Value = 7
$ws.Range("AL28").Value = 3
$ws.Range("AM28").Value = 10
$ws.Range("AN28").Value = 5
$ws.Range("AO28").Value = 30
$ws.Range("AP28").Value = 26
$ws.Range("AU28").Value = 2
$ws.Range("AY28").Value = 19
$ws.Range("BF28").Value = "2012-02-02"

# Row 29
$ws.Range("AD29").Value = 3
$ws.Range("AE29").Value = 24
$ws.Range("AG29").Value = 25
$ws.Range("AR29").Value = 26
$ws.Range("AT29").Value = 21
$ws.Range("AX29").Value = 20
$ws.Range("AY29").Value = 9
$ws.Range("BA29").Value = 25
$ws.Range("BF29").Value = "2012-02-02"

# Row 30
$ws.Range("D30").Value = 20
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 0.6
$ws.Range("H30").Value = 48.8
$ws.Range("I30").Value = 36.7
$ws.Range("J30").Value = 81.3
$ws.Range("L30").Value = 3.9
$ws.Range("M30").Value = 13.2
$ws.Range("N30").Value = 0.293
$ws.Range("O30").Value = 20.1
$ws.Range("P30").Value = 27.2
$ws.Range("Q30").Value = 0.737
$ws.Range("R30").Value = 11.7
$ws.Range("S30").Value = 31.1
$ws.Range("U30").Value = 20.3
$ws.Range("W30").Value = 8.300000000000001
$ws.Range("X30").Value = 6.3
$ws.Range("Z30").Value = 22.1
$ws.Range("AB30").Value = 97.2
$ws.Range("AC30").Value = 1.1
$ws.Range("AD30").Value = 27
$ws.Range("AG30").Value = 11
$ws.Range("AI30").Value = 10
$ws.Range("AJ30").Value = 11
$ws.Range("AO30").Value = 6
$ws.Range("AP30").Value = 4
$ws.Range("AQ30").Value = 19
$ws.Range("AR30").Value = 11
$ws.Range("AW30").Value = 11
$ws.Range("AX30").Value = 3
$ws.Range("AY30").Value = 24
$ws.Range("BA30").Value = 7
$ws.Range("BB30").Value = 10
$ws.Range("BC30").Value = 17
$ws.Range("BF30").Value = "2012-02-02"

# Row 31
$ws.Range("AD31").Value = 8
$ws.Range("AQ31").Value = 24
$ws.Range("AR31").Value = 8
$ws.Range("AT31").Value = 19
$ws.Range("BF31").Value = "2012-02-02"
